# Update "Seat Assignments" sheet:
#  1. Add a BookingEmailSent timestamp (col F) for the contestant in row 2.
#  2. Re-add a previously-deleted contestant record (the one that was sitting
#     in row 6 - "47ec0204-7cf6-4afb-9b9e-a91324cdf4a9") at the end of its
#     RecordDay block (just before the row that already has BookingEmailSent
#     / ConfirmedRSVP timestamps), this time also stamping its own
#     BookingEmailSent timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seat Assignments")

# --- 1. Stamp BookingEmailSent for row 2 -----------------------------------
# Copy the date/time number format already used lower in the column (F44)
# so the new cell gets the same style (rather than minting a new one).
$ws.Range("F44").Copy()
$ws.Range("F2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F2").Value = 46000.53846564815

# --- 2. Move the row-6 record to the end of its block ----------------------
# Capture the data currently sitting in row 6 before it is removed.
$a6 = $ws.Range("A6").Value2
$b6 = $ws.Range("B6").Value2
$c6 = $ws.Range("C6").Value2
$d6 = $ws.Range("D6").Value2
$e6 = $ws.Range("E6").Value2

# Remove row 6 - everything below shifts up one row. After this, the record
# that used to be row 44 ("A1", which already has BookingEmailSent /
# ConfirmedRSVP stamps) is now at row 43.
$ws.Rows.Item(6).Delete()

# Insert a fresh row right before that record, making room at row 43 for
# the restored record.
$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value = $a6
$ws.Range("B43").Value = $b6
$ws.Range("C43").Value = $c6
$ws.Range("D43").Value = $d6
$ws.Range("E43").Value = $e6

# Stamp its BookingEmailSent timestamp too, matching the styling used by
# the neighbouring timestamp cells.
$ws.Range("F44").Copy()
$ws.Range("F43").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F43").Value = 46000.538291944446
